$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric (e.g. "1.00", "0.999") must be forced to
# text format first, otherwise Excel auto-converts them to numbers and mangles
# formatting (e.g. "1.00" -> 1). Source data is a scraped text column.
$textForcedCells = @(
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D11",
    "D13",
    "D15",
    "D19",
    "D20",
    "D21",
    "D22",
    "D24",
    "D25",
    "D26",
    "D28",
    "D29",
    "D30",
    "D32",
    "D33",
    "D34",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51",
)
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply all updated values (prices, 1h volume %, and a few coin name/link
# corrections further down the table).
$ws.Range("D2").Value = "87.066.29"
$ws.Range("E2").Value = "  +3.40%  "
$ws.Range("D3").Value = "3.257.72"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "212.06"
$ws.Range("E5").Value = "  -3.50%  "
$ws.Range("D6").Value = "624.68"
$ws.Range("E6").Value = "  -1.96%  "
$ws.Range("D7").Value = "0.362"
$ws.Range("E7").Value = "  +10.77%  "
$ws.Range("D8").Value = "0.676"
$ws.Range("E8").Value = "  +14.62%  "
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "3.251.08"
$ws.Range("E10").Value = "  -1.11%  "
$ws.Range("D11").Value = "0.571"
$ws.Range("E11").Value = "  -4.70%  "
$ws.Range("E12").Value = "  +7.47%  "
$ws.Range("D13").Value = "0.0000254"
$ws.Range("E13").Value = "  -9.15%  "
$ws.Range("D14").Value = "3.859.81"
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("D15").Value = "33.50"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("E16").Value = "  -2.49%  "
$ws.Range("D17").Value = "86.707.45"
$ws.Range("E17").Value = "  +3.00%  "
$ws.Range("D18").Value = "3.310.31"
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("D19").Value = "3.11"
$ws.Range("E19").Value = "  -2.68%  "
$ws.Range("D20").Value = "13.98"
$ws.Range("E20").Value = "  -3.59%  "
$ws.Range("D21").Value = "431.71"
$ws.Range("E21").Value = "  -4.23%  "
$ws.Range("D22").Value = "8.84"
$ws.Range("E22").Value = "  -3.62%  "
$ws.Range("E23").Value = "  +0.87%  "
$ws.Range("D24").Value = "7.27"
$ws.Range("E24").Value = "  -2.55%  "
$ws.Range("D25").Value = "12.45"
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("D26").Value = "5.09"
$ws.Range("E26").Value = "  -4.02%  "
$ws.Range("D27").Value = "3.423.74"
$ws.Range("E27").Value = "  -0.94%  "
$ws.Range("D28").Value = "75.54"
$ws.Range("E28").Value = "  -3.42%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0000128"
$ws.Range("E29").Value = "  +1.15%  "
$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("E31").Value = "  +12.03%  "
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").Value = "8.68"
$ws.Range("E33").Value = "  -5.87%  "
$ws.Range("D34").Value = "541.00"
$ws.Range("E34").Value = "  -5.23%  "
$ws.Range("E35").Value = "  -8.24%  "
$ws.Range("D36").Value = "1.93"
$ws.Range("E36").Value = "  -4.74%  "
$ws.Range("D37").Value = "6.99"
$ws.Range("E37").Value = "  +12.04%  "
$ws.Range("D38").Value = "0.136"
$ws.Range("E38").Value = "  -11.56%  "
$ws.Range("D39").Value = "22.27"
$ws.Range("E39").Value = "  -4.64%  "
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "21.64"
$ws.Range("E41").Value = "  +3.40%  "
$ws.Range("D42").Value = "1.99"
$ws.Range("E42").Value = "  -3.14%  "
$ws.Range("D43").Value = "0.390"
$ws.Range("E43").Value = "  -5.91%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "2.91"
$ws.Range("E44").Value = "  -4.90%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "154.12"
$ws.Range("E46").Value = "  -3.60%  "
$ws.Range("D47").Value = "178.31"
$ws.Range("E47").Value = "  -6.86%  "
$ws.Range("D48").Value = "44.71"
$ws.Range("E48").Value = "  -0.98%  "
$ws.Range("D49").Value = "1.28"
$ws.Range("E49").Value = "  -4.83%  "
$ws.Range("D50").Value = "4.19"
$ws.Range("E50").Value = "  -1.84%  "
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").Value = "0.620"
$ws.Range("E51").Value = "  -3.94%  "

# Restore the default ("Normal") cell style on the cells we temporarily forced
# to text; the underlying cell stays text-typed once a string has been
# committed to it, so this only affects display formatting, not the value.
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).Style = "Normal"
}
